$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 21.97750833333333
$ws.Cells.Item(2, 8).Value = 65.932525
$ws.Cells.Item(2, 9).Value = 0.5427578249542736
$ws.Cells.Item(2, 10).Value = 0.5427578249542736
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 31.82741333333333
$ws.Cells.Item(2, 14).Value = 95.48223999999999
$ws.Cells.Item(2, 15).Value = 0.114390792932228
$ws.Cells.Item(2, 16).Value = 0.114390792932228
$ws.Cells.Item(2, 17).Value = 699.4872417617777
$ws.Cells.Item(2, 18).Value = 6295.385175855999
$ws.Cells.Item(2, 19).Value = 0.06208649796669077
$ws.Cells.Item(2, 20).Value = 0.06208649796669078

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 21.97750833333333
$ws.Cells.Item(3, 8).Value = 65.932525
$ws.Cells.Item(3, 9).Value = 0.5427578249542736
$ws.Cells.Item(3, 10).Value = 0.5427578249542736
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 85.46317833333335
$ws.Cells.Item(3, 14).Value = 256.389535
$ws.Cells.Item(3, 15).Value = 0.307162904935779
$ws.Cells.Item(3, 16).Value = 0.307162904935779
$ws.Cells.Item(3, 17).Value = 1878.267714013986
$ws.Cells.Item(3, 18).Value = 16904.40942612588
$ws.Cells.Item(3, 19).Value = 0.1667150701895797
$ws.Cells.Item(3, 20).Value = 0.1667150701895797

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 21.97750833333333
$ws.Cells.Item(4, 8).Value = 65.932525
$ws.Cells.Item(4, 9).Value = 0.5427578249542736
$ws.Cells.Item(4, 10).Value = 0.5427578249542736
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 122.2478306666667
$ws.Cells.Item(4, 14).Value = 366.743492
$ws.Cells.Item(4, 15).Value = 0.4393704929064738
$ws.Cells.Item(4, 16).Value = 0.4393704929064738
$ws.Cells.Item(4, 17).Value = 2686.702717208589
$ws.Cells.Item(4, 18).Value = 24180.3244548773
$ws.Cells.Item(4, 19).Value = 0.2384717730790048
$ws.Cells.Item(4, 20).Value = 0.2384717730790048

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 21.97750833333333
$ws.Cells.Item(5, 8).Value = 65.932525
$ws.Cells.Item(5, 9).Value = 0.5427578249542736
$ws.Cells.Item(5, 10).Value = 0.5427578249542736
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 38.69562533333333
$ws.Cells.Item(5, 14).Value = 116.086876
$ws.Cells.Item(5, 15).Value = 0.1390758092255191
$ws.Cells.Item(5, 16).Value = 0.1390758092255191
$ws.Cells.Item(5, 17).Value = 850.4334282268777
$ws.Cells.Item(5, 18).Value = 7653.900854041899
$ws.Cells.Item(5, 19).Value = 0.07548448371899825
$ws.Cells.Item(5, 20).Value = 0.07548448371899825

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 12.688376
$ws.Cells.Item(6, 8).Value = 38.065128
$ws.Cells.Item(6, 9).Value = 0.3133528721960219
$ws.Cells.Item(6, 10).Value = 0.3133528721960219
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 31.82741333333333
$ws.Cells.Item(6, 14).Value = 95.48223999999999
$ws.Cells.Item(6, 15).Value = 0.114390792932228
$ws.Cells.Item(6, 16).Value = 0.114390792932228
$ws.Cells.Item(6, 17).Value = 403.8381874807466
$ws.Cells.Item(6, 18).Value = 3634.54368732672
$ws.Cells.Item(6, 19).Value = 0.03584468351809405
$ws.Cells.Item(6, 20).Value = 0.03584468351809406

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 12.688376
$ws.Cells.Item(7, 8).Value = 38.065128
$ws.Cells.Item(7, 9).Value = 0.3133528721960219
$ws.Cells.Item(7, 10).Value = 0.3133528721960219
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 85.46317833333335
$ws.Cells.Item(7, 14).Value = 256.389535
$ws.Cells.Item(7, 15).Value = 0.307162904935779
$ws.Cells.Item(7, 16).Value = 0.307162904935779
$ws.Cells.Item(7, 17).Value = 1084.388940848387
$ws.Cells.Item(7, 18).Value = 9759.500467635482
$ws.Cells.Item(7, 19).Value = 0.0962503784937
$ws.Cells.Item(7, 20).Value = 0.0962503784937

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 12.688376
$ws.Cells.Item(8, 8).Value = 38.065128
$ws.Cells.Item(8, 9).Value = 0.3133528721960219
$ws.Cells.Item(8, 10).Value = 0.3133528721960219
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 122.2478306666667
$ws.Cells.Item(8, 14).Value = 366.743492
$ws.Cells.Item(8, 15).Value = 0.4393704929064738
$ws.Cells.Item(8, 16).Value = 0.4393704929064738
$ws.Cells.Item(8, 17).Value = 1551.126440682997
$ws.Cells.Item(8, 18).Value = 13960.13796614698
$ws.Cells.Item(8, 19).Value = 0.1376780059104254
$ws.Cells.Item(8, 20).Value = 0.1376780059104254

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 12.688376
$ws.Cells.Item(9, 8).Value = 38.065128
$ws.Cells.Item(9, 9).Value = 0.3133528721960219
$ws.Cells.Item(9, 10).Value = 0.3133528721960219
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 38.69562533333333
$ws.Cells.Item(9, 14).Value = 116.086876
$ws.Cells.Item(9, 15).Value = 0.1390758092255191
$ws.Cells.Item(9, 16).Value = 0.1390758092255191
$ws.Cells.Item(9, 17).Value = 490.9846437844586
$ws.Cells.Item(9, 18).Value = 4418.861794060128
$ws.Cells.Item(9, 19).Value = 0.04357980427380242
$ws.Cells.Item(9, 20).Value = 0.04357980427380242

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 2.008189666666667
$ws.Cells.Item(10, 8).Value = 6.024569
$ws.Cells.Item(10, 9).Value = 0.04959436889042158
$ws.Cells.Item(10, 10).Value = 0.04959436889042158
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 31.82741333333333
$ws.Cells.Item(10, 14).Value = 95.48223999999999
$ws.Cells.Item(10, 15).Value = 0.114390792932228
$ws.Cells.Item(10, 16).Value = 0.114390792932228
$ws.Cells.Item(10, 17).Value = 63.91548257272888
$ws.Cells.Item(10, 18).Value = 575.23934315456
$ws.Cells.Item(10, 19).Value = 0.005673139182348745
$ws.Cells.Item(10, 20).Value = 0.005673139182348746

$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 2.008189666666667
$ws.Cells.Item(11, 8).Value = 6.024569
$ws.Cells.Item(11, 9).Value = 0.04959436889042158
$ws.Cells.Item(11, 10).Value = 0.04959436889042158
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 85.46317833333335
$ws.Cells.Item(11, 14).Value = 256.389535
$ws.Cells.Item(11, 15).Value = 0.307162904935779
$ws.Cells.Item(11, 16).Value = 0.307162904935779
$ws.Cells.Item(11, 17).Value = 171.6262716094906
$ws.Cells.Item(11, 18).Value = 1544.636444485415
$ws.Cells.Item(11, 19).Value = 0.01523355041683852
$ws.Cells.Item(11, 20).Value = 0.01523355041683852

$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 2.008189666666667
$ws.Cells.Item(12, 8).Value = 6.024569
$ws.Cells.Item(12, 9).Value = 0.04959436889042158
$ws.Cells.Item(12, 10).Value = 0.04959436889042158
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 122.2478306666667
$ws.Cells.Item(12, 14).Value = 366.743492
$ws.Cells.Item(12, 15).Value = 0.4393704929064738
$ws.Cells.Item(12, 16).Value = 0.4393704929064738
$ws.Cells.Item(12, 17).Value = 245.4968303172164
$ws.Cells.Item(12, 18).Value = 2209.471472854948
$ws.Cells.Item(12, 19).Value = 0.02179030230477002
$ws.Cells.Item(12, 20).Value = 0.02179030230477002

$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 2.008189666666667
$ws.Cells.Item(13, 8).Value = 6.024569
$ws.Cells.Item(13, 9).Value = 0.04959436889042158
$ws.Cells.Item(13, 10).Value = 0.04959436889042158
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 38.69562533333333
$ws.Cells.Item(13, 14).Value = 116.086876
$ws.Cells.Item(13, 15).Value = 0.1390758092255191
$ws.Cells.Item(13, 16).Value = 0.1390758092255191
$ws.Cells.Item(13, 17).Value = 77.70815493960488
$ws.Cells.Item(13, 18).Value = 699.3733944564439
$ws.Cells.Item(13, 19).Value = 0.006897376986464292
$ws.Cells.Item(13, 20).Value = 0.006897376986464292

$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 3.818218
$ws.Cells.Item(14, 8).Value = 11.454654
$ws.Cells.Item(14, 9).Value = 0.09429493395928291
$ws.Cells.Item(14, 10).Value = 0.09429493395928291
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 31.82741333333333
$ws.Cells.Item(14, 14).Value = 95.48223999999999
$ws.Cells.Item(14, 15).Value = 0.114390792932228
$ws.Cells.Item(14, 16).Value = 0.114390792932228
$ws.Cells.Item(14, 17).Value = 121.5240024827733
$ws.Cells.Item(14, 18).Value = 1093.71602234496
$ws.Cells.Item(14, 19).Value = 0.01078647226509445
$ws.Cells.Item(14, 20).Value = 0.01078647226509445

$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 3.818218
$ws.Cells.Item(15, 8).Value = 11.454654
$ws.Cells.Item(15, 9).Value = 0.09429493395928291
$ws.Cells.Item(15, 10).Value = 0.09429493395928291
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 85.46317833333335
$ws.Cells.Item(15, 14).Value = 256.389535
$ws.Cells.Item(15, 15).Value = 0.307162904935779
$ws.Cells.Item(15, 16).Value = 0.307162904935779
$ws.Cells.Item(15, 17).Value = 326.3170458495434
$ws.Cells.Item(15, 18).Value = 2936.85341264589
$ws.Cells.Item(15, 19).Value = 0.02896390583566078
$ws.Cells.Item(15, 20).Value = 0.02896390583566078

$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 3.818218
$ws.Cells.Item(16, 8).Value = 11.454654
$ws.Cells.Item(16, 9).Value = 0.09429493395928291
$ws.Cells.Item(16, 10).Value = 0.09429493395928291
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 122.2478306666667
$ws.Cells.Item(16, 14).Value = 366.743492
$ws.Cells.Item(16, 15).Value = 0.4393704929064738
$ws.Cells.Item(16, 16).Value = 0.4393704929064738
$ws.Cells.Item(16, 17).Value = 466.7688675124187
$ws.Cells.Item(16, 18).Value = 4200.919807611768
$ws.Cells.Item(16, 19).Value = 0.04143041161227353
$ws.Cells.Item(16, 20).Value = 0.04143041161227353

$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 3.818218
$ws.Cells.Item(17, 8).Value = 11.454654
$ws.Cells.Item(17, 9).Value = 0.09429493395928291
$ws.Cells.Item(17, 10).Value = 0.09429493395928291
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 38.69562533333333
$ws.Cells.Item(17, 14).Value = 116.086876
$ws.Cells.Item(17, 15).Value = 0.1390758092255191
$ws.Cells.Item(17, 16).Value = 0.1390758092255191
$ws.Cells.Item(17, 17).Value = 147.7483331689893
$ws.Cells.Item(17, 18).Value = 1329.734998520904
$ws.Cells.Item(17, 19).Value = 0.01311414424625416
$ws.Cells.Item(17, 20).Value = 0.01311414424625416
